# Allocation rule summary tables: add "Within 5 miles" and "Within 10 miles"
# of HFC production facility columns (F, G) to both the "Means" and
# "Standard Deviations" sheets, and refresh the recomputed statistics that
# shifted as a result (Total Cancer Risk / Total Respiratory rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Means"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Means")

# New column headers
$ws1.Cells.Item(1, 6).Value = "Within 5 miles of HFC production facility"
$ws1.Cells.Item(1, 7).Value = "Within 10 miles of HFC production facility"

# New column F (Within 5 miles) and G (Within 10 miles) values, rows 2-10
$ws1.Cells.Item(2, 6).Value = 62
$ws1.Cells.Item(2, 7).Value = 66

$ws1.Cells.Item(3, 6).Value = 36
$ws1.Cells.Item(3, 7).Value = 27

$ws1.Cells.Item(4, 6).Value = 3
$ws1.Cells.Item(4, 7).Value = 7.1

$ws1.Cells.Item(5, 6).Value = 2.9
$ws1.Cells.Item(5, 7).Value = 5.1

$ws1.Cells.Item(6, 6).Value = 80
$ws1.Cells.Item(6, 7).Value = 79

$ws1.Cells.Item(7, 6).Value = 2.8
$ws1.Cells.Item(7, 7).Value = 5.7

$ws1.Cells.Item(8, 6).Value = 5.5
$ws1.Cells.Item(8, 7).Value = 4.9

$ws1.Cells.Item(9, 6).Value = 120
$ws1.Cells.Item(9, 7).Value = 80

$ws1.Cells.Item(10, 6).Value = 0.53
$ws1.Cells.Item(10, 7).Value = 0.51

# Existing values in row 9 (Total Cancer Risk) and row 10 (Total Respiratory)
# that changed because of the new radii being included in the recompute
$ws1.Cells.Item(9, 2).Value = 26
$ws1.Cells.Item(9, 3).Value = 39
$ws1.Cells.Item(9, 4).Value = 110
$ws1.Cells.Item(9, 5).Value = 120

$ws1.Cells.Item(10, 2).Value = 0.32
$ws1.Cells.Item(10, 3).Value = 0.43
$ws1.Cells.Item(10, 4).Value = 0.52
$ws1.Cells.Item(10, 5).Value = 0.53

# ---------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# New column headers
$ws2.Cells.Item(1, 6).Value = "Within 5 mile of HFC production facility SD"
$ws2.Cells.Item(1, 7).Value = "Within 10 mile of HFC production facility SD"

# New column F (Within 5 miles SD) and G (Within 10 miles SD) values, rows 2-10
$ws2.Cells.Item(2, 6).Value = 26
$ws2.Cells.Item(2, 7).Value = 26

$ws2.Cells.Item(3, 6).Value = 26
$ws2.Cells.Item(3, 7).Value = 26

$ws2.Cells.Item(4, 6).Value = 2.7
$ws2.Cells.Item(4, 7).Value = 6.7

$ws2.Cells.Item(5, 6).Value = 8.5
$ws2.Cells.Item(5, 7).Value = 6.2

$ws2.Cells.Item(6, 6).Value = 33
$ws2.Cells.Item(6, 7).Value = 32

$ws2.Cells.Item(7, 6).Value = 2.3
$ws2.Cells.Item(7, 7).Value = 7

$ws2.Cells.Item(8, 6).Value = 9.6
$ws2.Cells.Item(8, 7).Value = 8.1

$ws2.Cells.Item(9, 6).Value = 57
$ws2.Cells.Item(9, 7).Value = 32

$ws2.Cells.Item(10, 6).Value = 0.085
$ws2.Cells.Item(10, 7).Value = 0.066

# Existing values in row 9 (Total Cancer Risk SD) and row 10
# (Total Respiratory SD) that changed because of the new radii being
# included in the recompute
$ws2.Cells.Item(9, 2).Value = 8.6
$ws2.Cells.Item(9, 3).Value = 24
$ws2.Cells.Item(9, 4).Value = 64
$ws2.Cells.Item(9, 5).Value = 61

$ws2.Cells.Item(10, 2).Value = 0.14
$ws2.Cells.Item(10, 3).Value = 0.084
$ws2.Cells.Item(10, 4).Value = 0.05
$ws2.Cells.Item(10, 5).Value = 0.052
